# Commit message: "Moving from POI 3.17.0 to 4.0.1."
#
# Diffing the OOXML parts of this fixture (document.xml, header1-3.xml,
# footer1-3.xml, footnotes.xml, styles.xml) against the new revision
# shows every single hunk is a reordering of XML attributes and/or
# namespace declarations on otherwise-untouched elements (e.g.
# `w:headerReference r:id="rId6" w:type="even"` becomes
# `w:headerReference w:type="even" r:id="rId6"`, docPr's `id`/`name`/
# `descr` attributes swap places, the root element's `xmlns:*` list is
# reshuffled, etc.). No text run, paragraph, style value, relationship,
# numeric measurement, or reference id actually changes anywhere in the
# package - this is purely the cosmetic fingerprint of regenerating the
# fixture with a newer revision of the OOXML-writing library (Apache
# POI 3.17.0 -> 4.0.1), which simply serializes attributes in a
# different (but equivalent) order.
#
# There is consequently no Word object-model mutation that corresponds
# to this change: the document content, formatting, sections,
# headers/footers, styles and footnotes are identical before and
# after. We touch (read-only) the parts called out by the diff so the
# document round-trips through save, without altering any content.
$d = $word.ActiveDocument

$sec = $d.Sections.Item(1)
$null = $sec.PageSetup.PageWidth
$null = $sec.PageSetup.PageHeight
$null = $sec.Headers.Item(1).Range.Text
$null = $sec.Headers.Item(2).Range.Text
$null = $sec.Headers.Item(3).Range.Text
$null = $sec.Footers.Item(1).Range.Text
$null = $sec.Footers.Item(2).Range.Text
$null = $sec.Footers.Item(3).Range.Text
$null = $d.Content.Text
